$d = $word.ActiveDocument

# The document's East Asian default font is being switched from
# "DejaVu Sans" to "Tahoma" across the style sheet. Apply this to every
# paragraph style that currently carries an explicit w:eastAsia value
# (Normal and Heading); styles that inherit it (Body Text, List, Caption,
# Index) are left alone for the East Asian slot.
$d.Styles("Normal").Font.NameFarEast  = "Tahoma"
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# List, Caption and Index gain an explicit complex-script font
# (w:rFonts/@w:cs="DejaVu Sans") that previously was only inherited from
# Normal/Body Text. Setting Font.NameBi writes that attribute.
$d.Styles("List").Font.NameBi    = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi   = "DejaVu Sans"
